$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.353
$ws.Range("B4").Value = 5.175

$ws.Range("A7").Value = -20.513

$ws.Range("D10").Value = -8.099000000000002

$ws.Range("B12").Value = 5.609

$ws.Range("D13").Value = -8

$ws.Range("A16").Value = -22.119

$ws.Range("B18").Value = 5.032

$ws.Range("B19").Value = 8.272

$ws.Range("B20").Value = 6.471000000000001

$ws.Range("A28").Value = -21.645

$ws.Range("A29").Value = -21.532

$ws.Range("D30").Value = -7.139

$ws.Range("B31").Value = 6.706999999999999

$ws.Range("A32").Value = -21.821

$ws.Range("A40").Value = -20.944
$ws.Range("B40").Value = 7.042
$ws.Range("D40").Value = -8.301

$ws.Range("B42").Value = 6.854000000000001

$ws.Range("D44").Value = -7.664999999999999

$ws.Range("B47").Value = 6.146

$ws.Range("B48").Value = 6.276

$ws.Range("A52").Value = -21.675

$ws.Range("A57").Value = -22.132

$ws.Range("B63").Value = 5.583

$ws.Range("B64").Value = 6.064

$ws.Range("A66").Value = -21.529

$ws.Range("B76").Value = 6.003

$ws.Range("B81").Value = 4.842000000000001

$ws.Range("B89").Value = 5.007
$ws.Range("D89").Value = -8.195000000000002

$ws.Range("D91").Value = -7.514

$ws.Range("B94").Value = 5.923

$ws.Range("A100").Value = -22.352
